# chartink_screener.xlsx -- "break out stock.yaml completed"
#
# Sheet "3 V 0.3" gains one more scraped row (row 5), and the bsecode
# value in row 4 (previously entered/left as text) becomes a proper
# number, matching how the rest of that column is stored.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3 V 0.3")

# --- Fix row 4: bsecode should be numeric, not text -------------------
$ws.Range("E4").Value = 532900

# --- Append row 5 with the next scraped screener snapshot --------------
$ws.Range("A5").Value = "12/06/2024 08:44:44"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "PAISALO"
$ws.Range("D5").Value = "Paisalo Digital Ltd"

# bsecode here stays textual (as scraped) -- force text storage via a
# leading apostrophe, then drop the quote-prefix style it introduces so
# the cell keeps the sheet's default (unstyled) formatting.
$ws.Range("E5").Value = "'532900"
$ws.Range("E5").Style = "Normal"

$ws.Range("F5").Value = 9.35
$ws.Range("G5").Value = 69.5
$ws.Range("H5").Value = 4297806
